$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tax_sched_single")

# 1. Insert a new column before column D ("child"), shifting zve_nokfb..R right by one.
$ws.Columns.Item(4).Insert()

# 2. Header + values for new "child" column (all FALSE)
$ws.Range("D1").Value = "child"
$ws.Range("D2").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $false
$ws.Range("D5").Value = $false
$ws.Range("D6").Value = $false

# 3. Fix up the abgst / abgst_tu calculation block.
#    After the column insert, the old Q/R formulas (abgst / abgst_tu) live in
#    R/S now and still reference the old columns; rewrite them plus add the
#    new "meanwages" style Q column.
$ws.Range("Q2").Formula = "=(H2>9000)*(H2<13996)*(997.8*(H2-9000)/10000+1400)*(H2-9000)/10000+(H2>13996)*(H2<54949)*((220.13*(H2-13996)/10000+2397)*(H2-13996)/10000+948.49)+(H2>54950)*(H2<260532)*(0.42*H2-8621.75)+(H2>260532)*(0.45*H2-16437.7)"
$ws.Range("Q2").Style = "Percent"
$ws.Range("Q2").NumberFormat = "0.00"

$ws.Range("R2").Formula = "=MAX((I2-801)*0.25,0)"
$ws.Range("S2").Formula = "=R2"

$ws.Range("Q3").Formula = "=(H3>T3)*(H3<U3)*(997.6*(H3-T3)/10000+1400)*(H3-T3)/10000+(H3>U3)*(H3<W3)*((228.74*(H3-U3)/10000+2397)*(H3-V3)/10000+971)+(H3>W3)*(H3<W3)*(0.42*H3-8239)+(H3>W3)*(0.45*H3-15761)"

$ws.Range("R3:R6").Formula = "=MAX((I3-801)*0.25,0)"
$ws.Range("S3:S6").Formula = "=R3"

$ws.Range("Q4").Formula = "=(H4>T4)*(H4<U4)*(912.17*(H4-T4)/10000+1400)*(H4-T4)/10000+(H4>U4)*(H4<V4)*((228.74*(H4-U4)/10000+2397)*(H4-U4)/10000+1038)+(H4>V4)*(H4<W4)*(0.42*H4-8172)+(H4>W4)*(0.45*H4-15694)"
$ws.Range("Q5").Formula = "=(H5>T5)*(H5<U5)*(883.74*(H5-T5)/10000+1500)*(H5-T5)/10000+(H5>U5)*(H5<V5)*((228.74*(H5-U5)/10000+2397)*(H5-U5)/10000+989)+(H5>V5)*(H5<W5)*(0.42*H5-7914)+(H5>W5)*(0.45*H5-15414)"
$ws.Range("Q6").Formula = "=(H6>9000)*(H6<13996)*(997.8*(H6-9000)/10000+1400)*(H6-9000)/10000+(H6>13996)*(H6<54949)*((220.13*(H6-13996)/10000+2397)*(H6-13996)/10000+948.49)+(H6>54950)*(H6<260532)*(0.42*H6-8621.75)+(H6>260532)*(0.45*H6-16437.7)"

$ws.Range("Q3:Q6").NumberFormat = "0.00"

# 4. Update selection to match the authored state (D7 was the new active cell).
$ws.Range("D7").Select()
